$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, [string]$text) {
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range('D2') '43.387.18'
Set-TextValue $ws.Range('E2') '  -1.26%  '
Set-TextValue $ws.Range('D3') '2.371.78'
Set-TextValue $ws.Range('E3') '  +4.78%  '
Set-TextValue $ws.Range('E4') '  -0.04%  '
Set-TextValue $ws.Range('D5') '234.02'
Set-TextValue $ws.Range('E5') '  +0.25%  '
Set-TextValue $ws.Range('D6') '0.650'
Set-TextValue $ws.Range('E6') '  -0.52%  '
Set-TextValue $ws.Range('D7') '72.49'
Set-TextValue $ws.Range('E7') '  +13.43%  '
Set-TextValue $ws.Range('E8') '  +0.00%  '
Set-TextValue $ws.Range('E9') '  +6.64%  '
Set-TextValue $ws.Range('D10') '0.0976'
Set-TextValue $ws.Range('E10') '  -0.04%  '
Set-TextValue $ws.Range('E11') '  -2.46%  '
Set-TextValue $ws.Range('D12') '27.38'
Set-TextValue $ws.Range('E12') '  +3.19%  '
Set-TextValue $ws.Range('D13') '2.729.78'
Set-TextValue $ws.Range('E13') '  +5.03%  '
Set-TextValue $ws.Range('E14') '  +0.46%  '
Set-TextValue $ws.Range('D15') '16.00'
Set-TextValue $ws.Range('E15') '  +2.23%  '
Set-TextValue $ws.Range('D16') '6.29'
Set-TextValue $ws.Range('E16') '  +2.14%  '
Set-TextValue $ws.Range('D17') '0.864'
Set-TextValue $ws.Range('E17') '  +2.49%  '
Set-TextValue $ws.Range('D18') '2.373.62'
Set-TextValue $ws.Range('E18') '  +5.02%  '
Set-TextValue $ws.Range('D19') '43.406.35'
Set-TextValue $ws.Range('E19') '  -1.00%  '
Set-TextValue $ws.Range('D20') '0.0₃0995'
Set-TextValue $ws.Range('E20') '  +1.50%  '
Set-TextValue $ws.Range('D21') '6.36'
Set-TextValue $ws.Range('E21') '  +2.90%  '
Set-TextValue $ws.Range('D22') '74.65'
Set-TextValue $ws.Range('E22') '  +1.11%  '
Set-TextValue $ws.Range('D23') '249.91'
Set-TextValue $ws.Range('E23') '  -0.03%  '
Set-TextValue $ws.Range('E24') '  -0.12%  '
Set-TextValue $ws.Range('D25') '3.71'
Set-TextValue $ws.Range('E25') '  +3.81%  '
Set-TextValue $ws.Range('E26') '  -0.02%  '
Set-TextValue $ws.Range('D27') '2.27'
Set-TextValue $ws.Range('E27') '  +1.15%  '
Set-TextValue $ws.Range('D28') '10.01'
Set-TextValue $ws.Range('E28') '  +1.02%  '
Set-TextValue $ws.Range('D29') '22.51'
Set-TextValue $ws.Range('E29') '  +2.82%  '
Set-TextValue $ws.Range('D30') '174.12'
Set-TextValue $ws.Range('E30') '  +0.01%  '
Set-TextValue $ws.Range('E31') '  +5.92%  '
Set-TextValue $ws.Range('E32') '  -5.99%  '
Set-TextValue $ws.Range('E33') '  +0.12%  '
Set-TextValue $ws.Range('D34') '5.00'
Set-TextValue $ws.Range('E34') '  +0.84%  '
Set-TextValue $ws.Range('D35') '0.0695'
Set-TextValue $ws.Range('E35') '  +1.21%  '
Set-TextValue $ws.Range('D36') '5.07'
Set-TextValue $ws.Range('E36') '  +2.13%  '
Set-TextValue $ws.Range('B37') 'THORChain'
Set-TextValue $ws.Range('C37') 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue $ws.Range('D37') '6.63'
Set-TextValue $ws.Range('E37') '  +3.01%  '
Set-TextValue $ws.Range('B38') 'LidoDAOToken'
Set-TextValue $ws.Range('C38') 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
Set-TextValue $ws.Range('D38') '2.46'
Set-TextValue $ws.Range('E38') '  +6.91%  '
Set-TextValue $ws.Range('E39') '  +0.00%  '
Set-TextValue $ws.Range('E40') '  -0.18%  '
Set-TextValue $ws.Range('E41') '  +0.08%  '
Set-TextValue $ws.Range('D42') '8.90'
Set-TextValue $ws.Range('E42') '  +1.48%  '
Set-TextValue $ws.Range('D43') '18.44'
Set-TextValue $ws.Range('E43') '  +5.70%  '
Set-TextValue $ws.Range('E44') '  +9.26%  '
Set-TextValue $ws.Range('D45') '100.21'
Set-TextValue $ws.Range('E45') '  +1.35%  '
Set-TextValue $ws.Range('E46') '  +1.20%  '
Set-TextValue $ws.Range('D47') '4.46'
Set-TextValue $ws.Range('E47') '  -1.77%  '
Set-TextValue $ws.Range('D48') '0.0956'
Set-TextValue $ws.Range('E48') '  +0.38%  '
Set-TextValue $ws.Range('D49') '1.443.35'
Set-TextValue $ws.Range('E49') '  -0.88%  '
Set-TextValue $ws.Range('D50') '2.600.29'
Set-TextValue $ws.Range('E50') '  +5.09%  '
Set-TextValue $ws.Range('E51') '  -4.14%  '
